# Applies the "Add files via upload" edit:
#   - Column C (user fields, rows 6-9) shift up by one: the former C6
#     ("product_id") is dropped, each row takes the value that used to sit
#     one row below it, and the trailing C10 cell is cleared.
#   - Column E (product fields, rows 13-17) shift down by one: a new
#     "product_image" string is inserted at E13, every following row takes
#     the value that used to sit one row above it, and a new E17 cell is
#     created to hold the value that fell off the end.
#   - The worksheet selection moves from F21 to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: shift values in C6:C10 up by one row, clearing C10 ---
# (snapshot the original values first so the sequential writes below don't
# clobber a value before it has been read)
$c7 = $ws.Range("C7").Value2
$c8 = $ws.Range("C8").Value2
$c9 = $ws.Range("C9").Value2
$c10 = $ws.Range("C10").Value2

$ws.Range("C6").Value2 = $c7
$ws.Range("C7").Value2 = $c8
$ws.Range("C8").Value2 = $c9
$ws.Range("C9").Value2 = $c10
$ws.Range("C10").ClearContents()

# --- Column E: shift values in E13:E17 down by one row, inserting the new value at E13 ---
$e13 = $ws.Range("E13").Value2
$e14 = $ws.Range("E14").Value2
$e15 = $ws.Range("E15").Value2
$e16 = $ws.Range("E16").Value2

$ws.Range("E17").Value2 = $e16
$ws.Range("E16").Value2 = $e15
$ws.Range("E15").Value2 = $e14
$ws.Range("E14").Value2 = $e13
$ws.Range("E13").Value2 = "product_image"

# --- Selection moves to E3 ---
$ws.Range("E3").Select()
